$wb = $excel.ActiveWorkbook

# --- kpis (sheet1): add row 2 "Wins / 2 / games / football" ---
$kpis = $wb.Worksheets.Item("kpis")
$kpis.Range("A2").Value = "Wins"
$kpis.Range("B2").Value = 2
$kpis.Range("C2").Value = "games"
$kpis.Range("D2").Value = "football"
$kpis.Range("D3").Select() | Out-Null

# --- oil (sheet4): add quarterly rows 2-5 ---
$oil = $wb.Worksheets.Item("oil")
$oil.Range("A2").Value = 45748
$oil.Range("A2").NumberFormat = "mmm-yy"
$oil.Range("B2").Value = 100
$oil.Range("C2").Value = 1.1

$oil.Range("A3").Value = 45778
$oil.Range("A3").NumberFormat = "mmm-yy"
$oil.Range("B3").Value = 110
$oil.Range("C3").Value = 1.5

$oil.Range("A4").Value = 45809
$oil.Range("A4").NumberFormat = "mmm-yy"
$oil.Range("B4").Value = 120
$oil.Range("C4").Value = 1.4

$oil.Range("A5").Value = 45839
$oil.Range("A5").NumberFormat = "mmm-yy"
$oil.Range("B5").Value = 140
$oil.Range("C5").Value = 1.7

$oil.Range("A2:C5").Select() | Out-Null

# --- gas (sheet5): add quarterly rows 2-5 ---
$gas = $wb.Worksheets.Item("gas")
$gas.Range("A2").Value = 45748
$gas.Range("A2").NumberFormat = "mmm-yy"
$gas.Range("B2").Value = 200
$gas.Range("C2").Value = 1.1

$gas.Range("A3").Value = 45778
$gas.Range("A3").NumberFormat = "mmm-yy"
$gas.Range("B3").Value = 244
$gas.Range("C3").Value = 1.4

$gas.Range("A4").Value = 45809
$gas.Range("A4").NumberFormat = "mmm-yy"
$gas.Range("B4").Value = 299
$gas.Range("C4").Value = 2.3

$gas.Range("A5").Value = 45839
$gas.Range("A5").NumberFormat = "mmm-yy"
$gas.Range("B5").Value = 140
$gas.Range("C5").Value = 4

$gas.Range("C6").Select() | Out-Null

# --- footnotes (sheet9): becomes the active tab ---
$footnotes = $wb.Worksheets.Item("footnotes")
$footnotes.Activate() | Out-Null
